# PowerShell data table (row -> values) generated from the row permutation in the diff
$rowData = @{
  2 = @{ D=44400; I='Primera'; J=120; K=9000; L=10000; M=9500; N='$/caja 60 unidades'; P=158; Q=60 }
  3 = @{ D=44827; I='Primera'; J=120; K=6000; L=7000; M=6500; N='$/caja 60 unidades'; P=108; Q=60 }
  4 = @{ D=44627; I='Primera'; J=120; K=4000; L=4500; M=4250; N='$/caja 60 unidades'; P=71; Q=60 }
  5 = @{ D=44362; I='Primera'; J=120; K=8000; L=9000; M=8500; N='$/caja 60 unidades'; P=142; Q=60 }
  6 = @{ D=44494; I='Primera'; J=120; K=5000; L=6000; M=5500; N='$/caja 60 unidades'; P=92; Q=60 }
  7 = @{ D=44421; I='Primera'; J=100; K=8000; L=9000; M=8500; N='$/caja 60 unidades'; P=142; Q=60 }
  8 = @{ D=44764; I='Primera'; J=120; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; P=125; Q=60 }
  9 = @{ D=44785; I='Primera'; J=130; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; P=125; Q=60 }
  10 = @{ D=44281; I='Primera'; J=120; K=5500; L=6000; M=5750; N='$/caja 60 unidades'; P=96; Q=60 }
  11 = @{ D=44603; I='Primera'; J=140; K=5500; L=6000; M=5750; N='$/caja 60 unidades'; P=96; Q=60 }
  12 = @{ D=44740; I='Primera'; J=120; K=6000; L=7000; M=6500; N='$/caja 60 unidades'; P=108; Q=60 }
  13 = @{ D=44669; I='Primera'; J=130; K=4500; L=5000; M=4750; N='$/caja 60 unidades'; P=79; Q=60 }
  14 = @{ D=44935; I='Primera'; J=120; K=6000; L=7000; M=6500; N='$/caja 60 unidades'; P=108; Q=60 }
  15 = @{ D=44589; I='Primera'; J=110; K=5000; L=6000; M=5500; N='$/caja 60 unidades'; P=92; Q=60 }
  16 = @{ D=44242; I='Primera'; J=160; K=5000; L=5500; M=5250; N='$/caja 60 unidades'; P=88; Q=60 }
  17 = @{ D=44657; I='Primera'; J=100; K=5000; L=5500; M=5250; N='$/caja 60 unidades'; P=88; Q=60 }
  18 = @{ D=44967; I='Segunda'; J=50; K=4500; L=5000; M=4850; N='$/caja 90 unidades'; P=54; Q=90 }
  19 = @{ D=44963; I='Primera'; J=130; K=4000; L=4500; M=4250; N='$/caja 60 unidades'; P=71; Q=60 }
  20 = @{ D=44382; I='Primera'; J=160; K=7000; L=8000; M=7438; N='$/caja 60 unidades'; P=124; Q=60 }
  21 = @{ D=44676; I='Primera'; J=120; K=4000; L=4500; M=4250; N='$/caja 60 unidades'; P=71; Q=60 }
  22 = @{ D=44760; I='Primera'; J=130; K=7000; L=7500; M=7250; N='$/caja 60 unidades'; P=121; Q=60 }
  23 = @{ D=44648; I='Primera'; J=120; K=6500; L=7000; M=6750; N='$/caja 60 unidades'; P=112; Q=60 }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..23) {
    $data = $rowData[$r]
    $ws.Range("D$r").Value = $data.D
    $ws.Range("I$r").Value = $data.I
    $ws.Range("J$r").Value = $data.J
    $ws.Range("K$r").Value = $data.K
    $ws.Range("L$r").Value = $data.L
    $ws.Range("M$r").Value = $data.M
    $ws.Range("N$r").Value = $data.N
    $ws.Range("P$r").Value = $data.P
    $ws.Range("Q$r").Value = $data.Q
}

Write-Output "Done updating rows 2-23"
